$d = $word.ActiveDocument

# Locate the last paragraph of the list ("Installing react router...")
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Installing react router: npm i react-router-dom*") {
        $target = $p
    }
}

# Insert a new paragraph right after it; it inherits the ListParagraph
# style + numbering (numId) from the paragraph it was split off from.
$target.Range.InsertParagraphAfter()

$last = $d.Paragraphs.Last
$last.Range.Text = "When using Route in React must use link instead of <a> tag."
